$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New (updated) Cypher query text, with backtick-quoted multi-word aliases ---
$query = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN ['NCATS-COP01'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS `Case ID` , coalesce(s.clinical_study_designation,'') AS `Study Code` , coalesce(s.clinical_study_type,'') AS  `Study Type`, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  `Neutered Status`
'@

# --- Remove the "Url" column's look-and-feel from the cell that will become the new
#     query cell: clone B2's plain wrap-text style (no hyperlink font/underline) onto A2
#     before the hyperlink annotation and its now-unused "Hyperlink" cell style go away. ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# B2's content is moving to C2 (a plain, non-wrapped cell), so give B2 back the
# worksheet's default (un-styled) look now that A2 has its own copy of the format.
$ws.Range("B2").Style = "Normal"

$ws.Range("A2").Hyperlinks.Delete()
$wb.Styles.Item("Hyperlink").Delete()

# --- Remove the "Url" column (A) by writing the final 3-column layout directly: what
#     used to be query/dbExcel/WebExcel (B/C/D) becomes A/B/C. ---
$ws.Range("A1").Value = "query"
$ws.Range("B1").Value = "dbExcel"
$ws.Range("C1").Value = "WebExcel"

$ws.Range("A2").Value = $query
$ws.Range("B2").Value = "TC02_Canine_Filter_Study-NCATS_Neo4jData.xlsx"
$ws.Range("C2").Value = "TC02_Canine_Filter_Study-NCATS_WebData.xlsx"

# Drop the now-unused trailing column content (old WebExcel column D)
$ws.Range("D1:D2").ClearContents()

# --- Column widths for the new 3-column layout ---
$ws.Columns.Item(1).ColumnWidth = 75.81640625
$ws.Columns.Item(2).ColumnWidth = 70.26953125
$ws.Columns.Item(3).ColumnWidth = 28.54296875

# --- Row sized to fit the longer wrapped query text ---
$ws.Rows.Item(2).RowHeight = 188.5

# --- View / selection matches the final state ---
$ws.Range("A2").Select()
